$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet keeps a distinct number format on the last data row (e.g. a
# "date-only" look) versus the regular rows above it (which show a full
# date/time format). Before appending the new day's totals, remember the
# special "last row" format that currently sits on row 84, then restyle
# row 84 to match the regular rows since it will no longer be the last row.
$lastRowFormat = $ws.Range("A84").NumberFormat
$regularRowFormat = $ws.Range("A83").NumberFormat

$ws.Range("A84").NumberFormat = $regularRowFormat

# Append the new day's data as row 85.
$ws.Range("A85").Value = 45672
$ws.Range("B85").Value = 201
$ws.Range("C85").Value = 199
$ws.Range("D85").Value = 197

# Row 85 is now the last row, so it gets the special "last row" format that
# used to belong to row 84.
$ws.Range("A85").NumberFormat = $lastRowFormat
